# Apply "Varias Asignaturas y Planes" changes:
#  - Existing rows 2-9 (Arquitectura de Computadoras / Primero) change their
#    Docente from "admin" to "Ing. Carlos Guzman" and their Semestre from
#    "Primero" to "Tercero".
#  - A new block of rows (10-17) is appended, duplicating the same Carrera /
#    Actividad / Tema / Trabajo Independiente pattern, but for the subject
#    "Arquitectura para Redes", Semestre "Segundo".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the teacher for the existing 8 rows.
$ws.Range("A2:A9").Value = "Ing. Carlos Guzman"

# Update the semester for the existing 8 rows.
$ws.Range("D2:D9").Value = "Tercero"

# Activity/topic pairs shared by both subject blocks (columns F and G).
$topics = @(
  @("Encuadre", "Firmar Encuadre"),
  @("Prueba Diagnóstico", "SN"),
  @("Introducción", "SN"),
  @("Contenidos Varios", "Tipos de contenidos"),
  @("Mantenimientos", "SN"),
  @("Reparación", "Manual de reparación"),
  @("Evaluación Unidad", "Subir Portafolio"),
  @("Evaluación", "SN")
)

# Append the new "Arquitectura para Redes" / "Segundo" block in rows 10-17.
for ($i = 0; $i -lt 8; $i++) {
  $r = 10 + $i
  $ws.Cells.Item($r, 1).Value = "Ing. Carlos Guzman"
  $ws.Cells.Item($r, 2).Value = "Arquitectura para Redes"
  $ws.Cells.Item($r, 3).Value = "Redes y Telecomunicaciones"
  $ws.Cells.Item($r, 4).Value = "Segundo"
  $ws.Cells.Item($r, 5).Value = $i + 1
  $ws.Cells.Item($r, 6).Value = $topics[$i][0]
  $ws.Cells.Item($r, 7).Value = $topics[$i][1]
}

# Mirror the author's leftover selection state from the edit.
$ws.Range("D2:D9").Select()
